$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds the same "last changed" date serial for every
# data row (2-158). Bump it from 46060 (2026-02-07) to 46061 (2026-02-08).
for ($row = 2; $row -le 158; $row++) {
    $ws.Range("C$row").Value = 46061
}
